$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for fdcad896... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-05 00:51:42"

# Sheet "zh-cn": Correspond Handoff Datetime / Correspond Handback DateTime for fdcad896... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-05 00:51:38"
$wsZhCn.Range("K4").Value = "2016-09-05 00:51:56"

# Sheet "de-de": Correspond Handback DateTime for fdcad896... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-09-05 00:52:08"
